# Oyvinds verson som funker med 4 grupper
# Applies the "9groups -> 4 groups-capable" update to the Parameters sheet:
#  - bumps the 32/28 day-length table (I2:AJ2 / I3:AJ3) to 60/49 and 11/6
#  - appends a mirrored block (CF:DG) on rows 2/3 with the new values and on
#    rows 5/6 with the old values (kept as a reference/backup of the 9-group numbers)
#  - stamps an (initially empty) "CD" marker column down through row 61
#  - widens the sheet with two new column-width bands and moves the selection

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Parameters")

function Set-RowValues($sheet, $addr, $vals) {
    $n = $vals.Count
    $arr = New-Object 'object[,]' 1, $n
    for ($i = 0; $i -lt $n; $i++) { $arr[0, $i] = $vals[$i] }
    $sheet.Range($addr).Value = $arr
}

# ---- repeating 4-group patterns (7 values x 4 repeats = 28 columns, I:AJ / CF:DG) ----
$newRow2 = @(60,60,60,60,49,49,49, 60,60,60,60,49,49,49, 60,60,60,60,49,49,49, 60,60,60,60,49,49,49)
$newRow3 = @(11,11,11,11,6,6,6, 11,11,11,11,6,6,6, 11,11,11,11,6,6,6, 11,11,11,11,6,6,6)
$oldRow2 = @(32,32,32,32,28,28,28, 32,32,32,32,28,28,28, 32,32,32,32,28,28,28, 32,32,32,32,28,28,28)
$oldRow3 = @(6.358,6.358,6.358,6.358,3.468,3.468,3.468, 6.358,6.358,6.358,6.358,3.468,3.468,3.468, 6.358,6.358,6.358,6.358,3.468,3.468,3.468, 6.358,6.358,6.358,6.358,3.468,3.468,3.468)

# 1) Update the existing I2:AJ2 / I3:AJ3 values in place (32/28 -> 60/49, 6.358/3.468 -> 11/6)
Set-RowValues $ws "I2:AJ2" $newRow2
Set-RowValues $ws "I3:AJ3" $newRow3

# 2) New mirrored block in CF:DG -- rows 2/3 carry the new values, rows 5/6 keep the old ones
Set-RowValues $ws "CF2:DG2" $newRow2
Set-RowValues $ws "CF3:DG3" $newRow3
Set-RowValues $ws "CF5:DG5" $oldRow2
Set-RowValues $ws "CF6:DG6" $oldRow3

# 3) CD marker column, rows 2-61, formatted like the bold "G1-style" header cells (s=3) but empty
$ws.Range("G1").Copy()
$ws.Range("CD2:CD61").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# 4) Make sure the farthest cell is touched so dimension/used-range reaches DG61
$ws.Range("DG61").Value = $ws.Range("DG61").Value()

# 5) New column-width bands: CE (col 83) wide, CF:EI (cols 84-139) narrow
#    (ColumnWidth is quantized to 1/6-character steps by this host, so these are
#    the closest achievable inputs to the 10.33203125 / 3.33203125 target widths)
$ws.Columns.Item(83).ColumnWidth = 9.5
$ws.Range($ws.Cells.Item(1,84), $ws.Cells.Item(1,139)).EntireColumn.ColumnWidth = 2.5

# 6) Sheet view: drop the frozen "topLeftCell" pin and move the selection to X22
$ws.Activate()
$ws.Range("X22").Select() | Out-Null

# 7) Best-effort: nudge the workbook window position (host may not persist this)
try {
    $excel.ActiveWindow.Left = 1020
    $excel.ActiveWindow.Top = 500
} catch {}
